$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal text value without Excel auto-converting
# numeric-looking strings into actual numbers (so type/format matches source).
function Set-TextValue($cellRef, $text, $formatDonor) {
    $target = $ws.Range($cellRef)
    $target.NumberFormat = "@"
    $target.Value = $text
    $target.Style = $ws.Range($formatDonor).Style
}

$ws.Range("D2").Value = "29.965.04"
$ws.Range("D3").Value = "1.940.98"
$ws.Range("E3").Value = "  +1.11%  "
Set-TextValue "D4" "1.009" "D6"
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue "D5" "334.82" "D6"
$ws.Range("E5").Value = "  +2.69%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  +0.35%  "
Set-TextValue "D8" "0.4141" "D6"
$ws.Range("E8").Value = "  +1.50%  "
Set-TextValue "D9" "0.08210" "D6"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").Value = "  -0.44%  "
Set-TextValue "D11" "23.86" "D6"
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("D12").Value = "1.963.96"
$ws.Range("E12").Value = "  +2.54%  "
Set-TextValue "D13" "6.106" "D6"
$ws.Range("E13").Value = "  +1.00%  "
Set-TextValue "D14" "7.316" "D6"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("E15").Value = "  +0.23%  "
Set-TextValue "D16" "0.06858" "D6"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("E18").Value = "  +0.06%  "
Set-TextValue "D19" "17.89" "D6"
$ws.Range("E19").Value = "  +0.53%  "
Set-TextValue "D20" "1.009" "D6"
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "29.960.70"
$ws.Range("E21").Value = "  +1.57%  "
Set-TextValue "D22" "5.648" "D6"
$ws.Range("E22").Value = "  +0.17%  "
Set-TextValue "D23" "11.93" "D6"
$ws.Range("E23").Value = "  +1.16%  "
Set-TextValue "D24" "2.194" "D6"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "2.196.59"
$ws.Range("E25").Value = "  +2.08%  "
Set-TextValue "D26" "6.605" "D6"
$ws.Range("E26").Value = "  -0.24%  "
Set-TextValue "D27" "157.23" "D6"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("E28").Value = "  +0.07%  "
Set-TextValue "D29" "2.106" "D6"
$ws.Range("E29").Value = "  -0.44%  "
Set-TextValue "D30" "121.39" "D6"
$ws.Range("E30").Value = "  +0.85%  "
Set-TextValue "D31" "1.017" "D6"
$ws.Range("E31").Value = "  -0.53%  "
Set-TextValue "D32" "0.09639" "D6"
$ws.Range("E32").Value = "  +0.73%  "
Set-TextValue "D33" "5.625" "D6"
$ws.Range("E33").Value = "  +1.76%  "
Set-TextValue "D34" "1.422" "D6"
$ws.Range("E34").Value = "  +2.85%  "
Set-TextValue "D35" "3.559" "D6"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("E36").Value = "  +6.87%  "
Set-TextValue "D37" "0.02295" "D6"
$ws.Range("E37").Value = "  +0.50%  "
Set-TextValue "D38" "1.221" "D6"
$ws.Range("E38").Value = "  +3.15%  "
Set-TextValue "D39" "0.5991" "D6"
$ws.Range("E39").Value = "  +0.17%  "
Set-TextValue "D40" "8.024" "D6"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("E41").Value = "  -0.59%  "
Set-TextValue "D42" "2.540" "D6"
$ws.Range("E42").Value = "  +5.78%  "
Set-TextValue "D43" "0.1857" "D6"
$ws.Range("E43").Value = "  +0.06%  "
Set-TextValue "D44" "1.279" "D6"
$ws.Range("E44").Value = "  -0.14%  "
Set-TextValue "D45" "12.46" "D6"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("E46").Value = "  -1.24%  "
Set-TextValue "D47" "0.5578" "D6"
$ws.Range("E47").Value = "  +0.06%  "
Set-TextValue "D48" "1.987" "D6"
$ws.Range("E48").Value = "  +1.42%  "
Set-TextValue "D49" "117.88" "D6"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D50" "72.98" "D6"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("B51").Value = "MXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D51" "2.429" "D6"
$ws.Range("E51").Value = "  -0.09%  "
